# Swap the "Enterprises (absolute #)" and "Enterprises density (per 1000 people)"
# rows so that density appears first (row 11) and absolute count second (row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Enterprises density (per 1000 people)"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.8"

$ws.Range("A12").Value = "Enterprises (absolute #)"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4200"
